$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.76'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '25.11'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.013'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.577'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.008'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8120'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8375'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1338'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03279'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06952'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02839'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09407'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001515'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0005970'
$ws.Range("E16").Value = '15OneONE'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006083'
$ws.Range("E17").Value = '16TigerCashTCH'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.498'
$ws.Range("E18").Value = '17LEOLEO'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.091'
$ws.Range("E19").Value = '18BTSETokenBTSE'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3197'
$ws.Range("E20").Value = '19BitpandaEcosystemTokenBEST'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1315'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.746'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04676'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1370'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001241'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004522'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009701'
$ws.Range("E27").Value = '26NitroExNTXBestin24h'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03660'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1355'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006241'
$ws.Range("E42").Value = '41KickTokenKICK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008077'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005287'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002039'
